$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.436.45"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3
$ws.Range("D3").Value = "3.489.78"
$ws.Range("E3").Value = "  -0.75%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.77"
$ws.Range("E5").Value = "  -1.18%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.23"
$ws.Range("E6").Value = "  +0.13%  "

# Row 7
$ws.Range("D7").Value = "3.487.02"
$ws.Range("E7").Value = "  -0.78%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("E9").Value = "  +0.84%  "

# Row 10
$ws.Range("E10").Value = "  +2.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.51"
$ws.Range("E11").Value = "  +6.16%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.429"
$ws.Range("E12").Value = "  +0.95%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000216"
$ws.Range("E13").Value = "  -2.03%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.21"
$ws.Range("E14").Value = "  +0.62%  "

# Row 15
$ws.Range("D15").Value = "4.080.50"
$ws.Range("E15").Value = "  -0.78%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.490.43"
$ws.Range("E16").Value = "  -0.84%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "67.296.39"
$ws.Range("E17").Value = "  -0.04%  "

# Row 18
$ws.Range("E18").Value = "  -0.68%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.51"
$ws.Range("E19").Value = "  +1.80%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.51"
$ws.Range("E20").Value = "  +1.82%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.80"
$ws.Range("E21").Value = "  +4.37%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "444.50"
$ws.Range("E22").Value = "  +0.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.628"
$ws.Range("E23").Value = "  +0.41%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.99"
$ws.Range("E24").Value = "  +0.80%  "

# Row 25
$ws.Range("D25").Value = "3.629.92"
$ws.Range("E25").Value = "  -0.80%  "

# Row 26
$ws.Range("E26").Value = "  -0.07%  "

# Row 27
$ws.Range("E27").Value = "  -2.91%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -3.01%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.64"
$ws.Range("E29").Value = "  +3.26%  "

# Row 30
$ws.Range("E30").Value = "  +0.17%  "

# Row 31
$ws.Range("E31").Value = "  +5.96%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.168"
$ws.Range("E32").Value = "  +2.75%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.14%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.57"
$ws.Range("E34").Value = "  -1.37%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.13"
$ws.Range("E35").Value = "  -0.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.85"
$ws.Range("E36").Value = "  +0.52%  "

# Row 37
$ws.Range("D37").Value = "3.480.31"
$ws.Range("E37").Value = "  -0.83%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.96"
$ws.Range("E38").Value = "  -0.60%  "

# Row 39
$ws.Range("E39").Value = "  -0.01%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.30"
$ws.Range("E40").Value = "  +6.16%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "176.88"
$ws.Range("E42").Value = "  -0.89%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0885"
$ws.Range("E43").Value = "  +0.61%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.43"
$ws.Range("E44").Value = "  -0.18%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.887"
$ws.Range("E45").Value = "  +0.70%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.90"
$ws.Range("E46").Value = "  +4.62%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.13"
$ws.Range("E47").Value = "  +3.06%  "

# Row 48
$ws.Range("E48").Value = "  +2.61%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.53"
$ws.Range("E49").Value = "  -3.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.60"
$ws.Range("E50").Value = "  +0.22%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.988"
$ws.Range("E51").Value = "  -0.60%  "
